# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet and
#    populate it with the per-fund holding breakdown for 2022-Q1.
# 2. Update the "总计" (totals) sheet: insert a new first data row for
#    2022-Q1 (11 funds, 0.86 亿元) and push the existing quarters down.

$wb = $excel.ActiveWorkbook

# Helper: write a value that must stay TEXT even when it "looks" numeric
# (fund codes with leading zeros, decimal-looking percentages, etc).
# Plain $cell.Value = "4.03" gets auto-coerced to a number by the COM
# layer, so we go through the classic quote-prefix trick and then reset
# the style back to Normal so no stray number-format/quote-prefix style
# is left behind on the cell.
function Set-TextValue($cell, [string]$text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q1" worksheet, positioned before "总计"
# ---------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$templateSheet = $wb.Worksheets.Item($wb.Worksheets.Count - 1)  # "2021-Q4" - same column layout

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Bring over the header-row (B1:H1) and column-A formatting from the
# analogous "2021-Q4" sheet so the new sheet matches the look of its
# siblings (bold/bordered header, bold/bordered index column).
$templateSheet.Range("B1:H1").Copy() | Out-Null
$newSheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$templateSheet.Cells.Item(2, 1).Copy() | Out-Null
$newSheet.Range("A2:A12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Header row
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Data rows: idx, fund code, fund name, fund size, stock position, position %, held value, rank
$rows = @(
    @(0,  "001167", "金鹰科技创新股票",           "4.03", "94.55", "4.63", "0.1866", 5),
    @(1,  "210009", "金鹰核心资源混合",           "3.86", "94.96", "4.77", "0.1841", 5),
    @(2,  "162102", "金鹰中小盘精选混合",         "4.60", "76.52", "3.46", "0.1592", 10),
    @(3,  "001415", "信诚新锐回报灵活配置混合A", "9.07", "24.72", "0.76", "0.0689", 3),
    @(4,  "001402", "信诚新选回报灵活配置混合A", "8.37", "22.05", "0.62", "0.0519", 6),
    @(5,  "003234", "信诚至利灵活配置混合A",     "8.99", "22.05", "0.57", "0.0512", 8),
    @(6,  "004157", "信诚至诚灵活配置混合A",     "7.32", "22.71", "0.67", "0.0490", 5),
    @(7,  "002046", "信诚新锐回报灵活配置混合B", "5.19", "24.72", "0.76", "0.0394", 3),
    @(8,  "003235", "信诚至利灵活配置混合C",     "5.30", "22.05", "0.57", "0.0302", 8),
    @(9,  "002030", "信诚新选回报灵活配置混合B", "3.85", "22.05", "0.62", "0.0239", 6),
    @(10, "004158", "信诚至诚灵活配置混合B",     "2.18", "22.71", "0.67", "0.0146", 5)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]                      # A: index (number)
    Set-TextValue $newSheet.Cells.Item($r, 2) $row[1]                 # B: 基金代码
    Set-TextValue $newSheet.Cells.Item($r, 3) $row[2]                 # C: 基金名称
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[3]                 # D: 基金规模
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[4]                 # E: 股票总仓位
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[5]                 # F: 仓位占比
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[6]                 # G: 持有市值(亿元)
    $newSheet.Cells.Item($r, 8).Value = $row[7]                       # H: 仓位排名 (number)
    $r++
}

# re-apply col-A formatting: the above writes reset A2's style to Normal via
# the value assignment on the header copy step; make sure the index column
# keeps its bold/bordered look after the data was written.
$templateSheet.Cells.Item(2, 1).Copy() | Out-Null
$newSheet.Range("A2:A12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# ---------------------------------------------------------------------
# Step 2: update "总计" - insert the 2022-Q1 row at the top, push others down
# ---------------------------------------------------------------------

$ws = $totalSheet

# Shift existing rows down one at a time, bottom-up, copying values and
# (for column A, which carries the bold/bordered style) formats too.
for ($src = 4; $src -ge 2; $src--) {
    $dst = $src + 1
    $ws.Range("A$src`:D$src").Copy() | Out-Null
    $ws.Range("A$dst`:D$dst").PasteSpecial(-4163) | Out-Null   # xlPasteValues

    $ws.Cells.Item($src, 1).Copy() | Out-Null
    $ws.Cells.Item($dst, 1).PasteSpecial(-4122) | Out-Null     # xlPasteFormats
}

# Write the new first data row: 2022-Q1 / 11 funds / 0.86 亿元
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "2022-Q1"
$ws.Cells.Item(2, 3).Value = 11
$ws.Cells.Item(2, 4).Value = 0.86

Write-Host "2022-Q1 sheet added and 总计 updated"
